$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Price (D) and Volume(1h) (E) columns to match the refreshed
# coin quote snapshot. Each cell's NumberFormat is (re)applied as Text
# immediately before its value is written so Excel keeps storing the
# numeric-looking strings (e.g. "311.85", "0.68%") as text, matching the
# original inline-string cell type instead of re-parsing them as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "311.85"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.68%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "38.36"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.83%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.127"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.15%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08106"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.01%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.481"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "5.60%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.960"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.35%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "8.316"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.92%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9407"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.31%"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-6.32%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1954"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.07%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09009"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.30%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03488"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.35%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09708"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.16%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001409"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.63%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006049"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.05%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.569"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-8.55%"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-3.90%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3466"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.38%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.69%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.025"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "6.44%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2492"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.64%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04376"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.18%"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.34%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004733"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.37%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003907"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "212.02%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02198"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "3.22%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05240"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.70%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007630"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.67%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01034"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "5.80%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1388"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "1.83%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002038"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-4.49%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009123"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.78%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006620"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "3.36%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000753"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.25%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003019"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "17.53%"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "68.82%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002109"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.25%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002008"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.25%"
